$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1985135.4
$ws.Range("I15").Value = 1985135.4
$ws.Range("K15").Value = 5955406.199999999
$ws.Range("M15").Value = -5955237.199999999
$ws.Range("H43").Value = 8028.8
$ws.Range("I43").Value = 4445
$ws.Range("J43").Value = 8924.75
$ws.Range("K43").Value = 4445
$ws.Range("L43").Value = 8924.75
$ws.Range("M43").Value = -4376
$ws.Range("N43").Value = -9062.75
$ws.Range("H51").Value = 39427.145
$ws.Range("I51").Value = 8995.5
$ws.Range("K51").Value = 8995.5
$ws.Range("M51").Value = -8511.5
$ws.Range("H110").Value = 64495
$ws.Range("J110").Value = 64495
$ws.Range("L110").Value = 64495
$ws.Range("N110").Value = -72675
$ws.Range("H113").Value = 5848.6
$ws.Range("I113").Value = 5590.125
$ws.Range("J113").Value = 6308.1113
$ws.Range("K113").Value = 5590.125
$ws.Range("L113").Value = 6308.1113
$ws.Range("M113").Value = -2336.125
$ws.Range("N113").Value = -12816.1113
$ws.Range("H138").Value = 3492.23
$ws.Range("J138").Value = 3563.7778
$ws.Range("L138").Value = 10691.3334
$ws.Range("N138").Value = -20971.3334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22136.066
$ws.Range("I32").Value = 10188.022
$ws.Range("J32").Value = 53060.41
$ws.Range("K32").Value = 10188.022
$ws.Range("L32").Value = 53060.41
$ws.Range("M32").Value = -9901.022000000001
$ws.Range("N32").Value = -53634.41
$ws.Range("H74").Value = 8047.6
$ws.Range("I74").Value = 3203.0833
$ws.Range("K74").Value = 3203.0833
$ws.Range("M74").Value = -2329.0833
$ws.Range("H77").Value = 8047.6
$ws.Range("I77").Value = 3203.0833
$ws.Range("K77").Value = 16015.4165
$ws.Range("M77").Value = -11647.4165
$ws.Range("H88").Value = 3997.7778
$ws.Range("I88").Value = 3496
$ws.Range("J88").Value = 4625
$ws.Range("K88").Value = 3496
$ws.Range("L88").Value = 4625
$ws.Range("M88").Value = -3090
$ws.Range("N88").Value = -5437
$ws.Range("H91").Value = 3997.7778
$ws.Range("I91").Value = 3496
$ws.Range("J91").Value = 4625
$ws.Range("K91").Value = 3496
$ws.Range("L91").Value = 4625
$ws.Range("M91").Value = -2092
$ws.Range("N91").Value = -7433
$ws.Range("H94").Value = 33197.6
$ws.Range("J94").Value = 33197.6
$ws.Range("L94").Value = 33197.6
$ws.Range("N94").Value = -34999.6
$ws.Range("H122").Value = 3404.75
$ws.Range("I122").Value = 2417.818
$ws.Range("K122").Value = 7253.454000000001
$ws.Range("M122").Value = -4803.454000000001
$ws.Range("H132").Value = 7775.65
$ws.Range("I132").Value = 4774.857
$ws.Range("K132").Value = 14324.571
$ws.Range("M132").Value = -11794.571

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 42033.582
$ws.Range("I82").Value = 22306.428
$ws.Range("K82").Value = 22306.428
$ws.Range("M82").Value = -21923.428
$ws.Range("H85").Value = 42033.582
$ws.Range("I85").Value = 22306.428
$ws.Range("K85").Value = 22306.428
$ws.Range("M85").Value = -20980.428
$ws.Range("H99").Value = 7622.048
$ws.Range("I99").Value = 11130.091
$ws.Range("K99").Value = 11130.091
$ws.Range("M99").Value = -9632.091
$ws.Range("H107").Value = 2244.3333
$ws.Range("I107").Value = 1628.8334
$ws.Range("K107").Value = 1628.8334
$ws.Range("M107").Value = 291.1666
$ws.Range("H134").Value = 3665948
$ws.Range("I134").Value = 1847.3226
$ws.Range("J134").Value = 17864338
$ws.Range("K134").Value = 5541.9678
$ws.Range("L134").Value = 53593014
$ws.Range("M134").Value = -3006.9678
$ws.Range("N134").Value = -53598084

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 979392.0600000001
$ws.Range("I31").Value = 14288.692
$ws.Range("J31").Value = 2119968.8
$ws.Range("K31").Value = 14288.692
$ws.Range("L31").Value = 2119968.8
$ws.Range("M31").Value = -13993.692
$ws.Range("N31").Value = -2120558.8
$ws.Range("H34").Value = 979392.0600000001
$ws.Range("I34").Value = 14288.692
$ws.Range("J34").Value = 2119968.8
$ws.Range("K34").Value = 14288.692
$ws.Range("L34").Value = 2119968.8
$ws.Range("M34").Value = -14086.692
$ws.Range("N34").Value = -2120372.8
$ws.Range("H69").Value = 24878.715
$ws.Range("I69").Value = 12872.4
$ws.Range("J69").Value = 54894.5
$ws.Range("K69").Value = 12872.4
$ws.Range("L69").Value = 54894.5
$ws.Range("M69").Value = -12123.4
$ws.Range("N69").Value = -56392.5
$ws.Range("H72").Value = 24878.715
$ws.Range("I72").Value = 12872.4
$ws.Range("J72").Value = 54894.5
$ws.Range("K72").Value = 38617.2
$ws.Range("L72").Value = 164683.5
$ws.Range("M72").Value = -34873.2
$ws.Range("N72").Value = -172171.5
$ws.Range("H122").Value = 2283.1667
$ws.Range("I122").Value = 2243.7693
$ws.Range("J122").Value = 2385.6
$ws.Range("K122").Value = 6731.3079
$ws.Range("L122").Value = 7156.799999999999
$ws.Range("M122").Value = -4281.3079
$ws.Range("N122").Value = -12056.8
$ws.Range("H132").Value = 7106
$ws.Range("I132").Value = 2906.4614
$ws.Range("K132").Value = 8719.3842
$ws.Range("M132").Value = -6189.3842

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 157.84314
$ws.Range("I2").Value = 77.61905
$ws.Range("J2").Value = 214
$ws.Range("K2").Value = 465.7143
$ws.Range("L2").Value = 1284
$ws.Range("M2").Value = -352.7143
$ws.Range("N2").Value = -1510
$ws.Range("H113").Value = 1169.12
$ws.Range("I113").Value = 709.2
$ws.Range("K113").Value = 2127.6
$ws.Range("M113").Value = 42.39999999999964
$ws.Range("H134").Value = 8530.235000000001
$ws.Range("I134").Value = 2377.625
$ws.Range("K134").Value = 7132.875
$ws.Range("M134").Value = -2062.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 308.375
$ws.Range("I2").Value = 264.3846
$ws.Range("J2").Value = 499
$ws.Range("K2").Value = 264.3846
$ws.Range("L2").Value = 499
$ws.Range("M2").Value = -151.3846
$ws.Range("N2").Value = -725
$ws.Range("H40").Value = 20000
$ws.Range("J40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("N40").Value = -20302
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H93").Value = 48358.168
$ws.Range("J93").Value = 48358.168
$ws.Range("L93").Value = 48358.168
$ws.Range("N93").Value = -52102.168
$ws.Range("H122").Value = 4196
$ws.Range("I122").Value = 4377.5
$ws.Range("J122").Value = 3978.2
$ws.Range("K122").Value = 13132.5
$ws.Range("L122").Value = 11934.6
$ws.Range("M122").Value = -10682.5
$ws.Range("N122").Value = -16834.6
$ws.Range("H132").Value = 5833.7646
$ws.Range("I132").Value = 5761.7334
$ws.Range("J132").Value = 6374
$ws.Range("K132").Value = 17285.2002
$ws.Range("L132").Value = 19122
$ws.Range("M132").Value = -14755.2002
$ws.Range("N132").Value = -24182

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 400828.56
$ws.Range("I7").Value = 7950.5
$ws.Range("J7").Value = 924666
$ws.Range("K7").Value = 7950.5
$ws.Range("L7").Value = 924666
$ws.Range("M7").Value = -7838.5
$ws.Range("N7").Value = -924890
$ws.Range("H40").Value = 6984.9614
$ws.Range("I40").Value = 6206.0557
$ws.Range("K40").Value = 6206.0557
$ws.Range("M40").Value = -6070.0557
$ws.Range("H93").Value = 1842.0344
$ws.Range("I93").Value = 1590
$ws.Range("K93").Value = 1590
$ws.Range("M93").Value = -342
$ws.Range("H100").Value = 2852.625
$ws.Range("I100").Value = 2428.5
$ws.Range("K100").Value = 2428.5
$ws.Range("M100").Value = -1887.5
$ws.Range("H126").Value = 400828.56
$ws.Range("I126").Value = 7950.5
$ws.Range("J126").Value = 924666
$ws.Range("K126").Value = 23851.5
$ws.Range("L126").Value = 2773998
$ws.Range("M126").Value = -21381.5
$ws.Range("N126").Value = -2778938
$ws.Range("H132").Value = 6535.0713
$ws.Range("I132").Value = 7744.636
$ws.Range("K132").Value = 23233.908
$ws.Range("M132").Value = -20703.908

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 65247.5
$ws.Range("J16").Value = 71996.664
$ws.Range("L16").Value = 71996.664
$ws.Range("N16").Value = -72580.664
$ws.Range("H62").Value = 20015500
$ws.Range("J62").Value = 20015500
$ws.Range("L62").Value = 20015500
$ws.Range("N62").Value = -20016748
$ws.Range("H65").Value = 20015500
$ws.Range("J65").Value = 20015500
$ws.Range("L65").Value = 100077500
$ws.Range("N65").Value = -100083740
$ws.Range("H132").Value = 2035109.4
$ws.Range("I132").Value = 2249.818
$ws.Range("J132").Value = 10420655
$ws.Range("K132").Value = 6749.454000000001
$ws.Range("L132").Value = 31261965
$ws.Range("M132").Value = -4219.454000000001
$ws.Range("N132").Value = -31267025
